{"js": "// \"update with this year stats\"\n// 1) SAFE's projected headcount moves from 60 to 75 students.\n// 2) Trim the \"(and I'm paid to do)\" aside down to a simple comma.\n\nconst headcountResults = context.document.body.search(\n  \"there will be over 60 students\",\n  { matchCase: true, matchWholeWord: false }\n);\nheadcountResults.load(\"items,text\");\nawait context.sync();\n\nif (headcountResults.items.length > 0) {\n  headcountResults.items[0].insertText(\n    \"there will be over 75 students\",\n    Word.InsertLocation.replace\n  );\n} else {\n  // Fallback: just the bare number, in case the surrounding wording differs.\n  const fallback = context.document.body.search(\"60\", { matchCase: true, matchWholeWord: true });\n  fallback.load(\"items\");\n  await context.sync();\n  if (fallback.items.length > 0) {\n    fallback.items[0].insertText(\"75\", Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n\nconst asideResults = context.document.body.search(\n  \"to do (and I\\u2019m paid to do) but\",\n  { matchCase: true, matchWholeWord: false }\n);\nasideResults.load(\"items,text\");\nawait context.sync();\n\nif (asideResults.items.length > 0) {\n  asideResults.items[0].insertText(\"to do, but\", Word.InsertLocation.replace);\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$rightApos = [char]8217\n\n# \"update with this year stats\"\n# 1) SAFE's projected headcount moves from 60 to 75 students.\n$find = $d.Content.Find\n$find.ClearFormatting()\n$find.Forward = $true\n$find.Wrap = 1\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n$find.MatchSoundsLike = $false\n$find.MatchAllWordForms = $false\n$find.Text = \"there will be over 60 students\"\n$find.Replacement.ClearFormatting()\n$find.Replacement.Text = \"there will be over 75 students\"\n$found = $find.Execute($find.Text, $find.MatchCase, $find.MatchWholeWord, $find.MatchWildcards, $find.MatchSoundsLike, $find.MatchAllWordForms, $find.Forward, $find.Wrap, $find.Format, $find.Replacement.Text, 2)\n\nif (-not $found) {\n    # Fallback: just the bare whole-word number, in case the surrounding wording differs.\n    $find3 = $d.Content.Find\n    $find3.ClearFormatting()\n    $find3.Forward = $true\n    $find3.Wrap = 1\n    $find3.Format = $false\n    $find3.MatchCase = $true\n    $find3.MatchWholeWord = $true\n    $find3.MatchWildcards = $false\n    $find3.Text = \"60\"\n    $find3.Replacement.ClearFormatting()\n    $find3.Replacement.Text = \"75\"\n    $find3.Execute($find3.Text, $find3.MatchCase, $find3.MatchWholeWord, $find3.MatchWildcards, $false, $false, $find3.Forward, $find3.Wrap, $find3.Format, $find3.Replacement.Text, 2)\n}\n\n# 2) Trim the \"(and I'm paid to do)\" aside down to a simple comma.\n$find2 = $d.Content.Find\n$find2.ClearFormatting()\n$find2.Forward = $true\n$find2.Wrap = 1\n$find2.Format = $false\n$find2.MatchCase = $true\n$find2.MatchWholeWord = $false\n$find2.MatchWildcards = $false\n$find2.MatchSoundsLike = $false\n$find2.MatchAllWordForms = $false\n$find2.Text = \"to do (and I\" + $rightApos + \"m paid to do) but\"\n$find2.Replacement.ClearFormatting()\n$find2.Replacement.Text = \"to do, but\"\n$find2.Execute($find2.Text, $find2.MatchCase, $find2.MatchWholeWord, $find2.MatchWildcards, $find2.MatchSoundsLike, $find2.MatchAllWordForms, $find2.Forward, $find2.Wrap, $find2.Format, $find2.Replacement.Text, 2)\n"}
